$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vcsa")

# Insert two new columns before column B (shifts old B.. to D..)
$ws.Range("B:C").Insert()

$ws.Range("B1").Value = "Config?"
$ws.Range("C1").Value = "Certs?"

$ws.Range("B2:C5").Value = $true

$ws.Range("A2").Select()
